$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.49%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.47%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.126"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.47%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07616"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.07%"

$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.279"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.74%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.616"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.87%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.457"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.82%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9089"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.41%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1278"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "28.41%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1808"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.31%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09106"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.21%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04293"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.63%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1044"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.71%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001263"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.13%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005785"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.34%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.350"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.59%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.46%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.956"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.26%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.06%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-4.83%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04046"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.86%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001272"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.68%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004039"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.58%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001273"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.02%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "24.78%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02426"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "0.42%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05233"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.72%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007830"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.05%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1297"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.56%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006809"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.63%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001900"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.44%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007361"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.03%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3343"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.31%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006900"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "8.32%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.22%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1077"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1,806.79%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.22%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.22%"
